# Security Vision Scope - update level-of-effort estimates, insert a new
# separator row before the "Investigate ..." block, and mark those three
# rows' Release column as unknown ("?").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Level of Effort (Days) bumps
$ws.Range("F11").Value = 20
$ws.Range("F16").Value = 40

# Insert a new blank separator row above row 28 ("Investigate SAML 2.0 ...")
# and give it the same formatting as the row above it (row 27), matching
# how the existing separator rows (14, 19, 31) look.
$ws.Rows(28).Insert()
$ws.Range("A27:H27").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)
$ws.Range("A28:H28").ClearContents()
$ws.Rows(28).RowHeight = $ws.Rows(27).RowHeight
$excel.CutCopyMode = $false

# The three rows that used to sit at 28-30 are now at 29-31; mark their
# Release value as unknown.
$ws.Range("A29").Value = "?"
$ws.Range("A30").Value = "?"
$ws.Range("A31").Value = "?"

# Move the active selection, matching where the author last clicked.
[void]$ws.Range("H18").Select()
